$ws = $excel.ActiveWorkbook.ActiveSheet

$ws.Range("D2").Value = "68.210.29"
$ws.Range("E2").Value = "  -0.64%  "
$ws.Range("D3").Value = "2.644.51"
$ws.Range("E3").Value = "  -0.24%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.46"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.35%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.43"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.18%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.542"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.66%  "
$ws.Range("E9").Value = "  +2.78%  "
$ws.Range("E10").Value = "  -1.40%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.26"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.59%  "
$ws.Range("E12").Value = "  +0.57%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "27.98"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.26%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000190"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.87%  "
$ws.Range("D15").Value = "3.125.18"
$ws.Range("E15").Value = "  -0.20%  "
$ws.Range("D16").Value = "68.212.94"
$ws.Range("E16").Value = "  -0.41%  "
$ws.Range("D17").Value = "2.654.48"
$ws.Range("E17").Value = "  +0.11%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.36"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.50%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "362.87"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.60%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.33"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.59%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.42"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +3.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.79"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.48%  "
$ws.Range("E23").Value = "  -2.79%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "75.32"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +3.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.74"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.53%  "
$ws.Range("E27").Value = "  +0.40%  "
$ws.Range("E28").Value = "  -1.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "556.86"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -3.15%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.05"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.52%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.40"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.47%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.86"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.62%  "
$ws.Range("E34").Value = "  +0.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.128"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.77%  "
$ws.Range("E36").Value = "  +0.15%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "161.89"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.75%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.62"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.92%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.372"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.06%  "
$ws.Range("E40").Value = "  -3.31%  "
$ws.Range("E41").Value = "  -1.53%  "
$ws.Range("B42").Value = "BabyDogeCoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D42").Value = "0.0₆0335"
$ws.Range("E42").Value = "  +4.76%  "
$ws.Range("B43").Value = "WhiteBITCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.79"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.28%  "
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "158.63"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.76%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.72"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.56%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "21.98"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.07%  "
$ws.Range("B49").Value = "Optimism"
$ws.Range("C49").Value = "https://coinranking.com/coin/n1p-s_gm1+optimism-op"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.68"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.92%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0780"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.04%  "
$ws.Range("E51").Value = "  -0.60%  "
